# Update the marksheet "Corr/total marks" figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking" -> Right column (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Row 12 "Total" -> Right column (B12): 30 -> 50, and the Max label (E12) reflects it: "30/84" -> "50/140"
$ws.Range("B12").Value = 50
$ws.Range("E12").Value = "50/140"
